$d = $word.ActiveDocument

# Locate the "Temas establecidos para la proxima reunion" heading, then
# find the (single, currently empty) list paragraph that follows it -
# that is where the new agenda items for the next meeting get listed.
$headingText = "Temas establecidos para la próxima reunión"
$count = $d.Paragraphs.Count
$anchorIndex = -1

for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd() -eq $headingText) {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find the 'Temas establecidos para la proxima reunion' heading"
}

$startIndex = $anchorIndex + 1

$items = @(
    "Gestión Judicial",
    "Ajuste de Cuotas",
    "Parametrización de Planes de Pagos",
    "Parametrización de Categoría de Contribuyentes",
    "Parametrización de Tipos de Contribuyentes",
    "Informe de Imputaciones de Pagos",
    "Prescripción",
    "Reversión de Cobros",
    "Administración de Zonas de Reparto"
)

for ($i = 0; $i -lt $items.Length; $i++) {
    $p = $d.Paragraphs.Item($startIndex + $i)
    $p.Range.InsertAfter($items[$i])
    if ($i -lt $items.Length - 1) {
        $p.Range.InsertParagraphAfter()
    }
}
